$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.2027491408934708
$ws.Range("C2").Value = 0.5360824742268041
$ws.Range("J2").Value = 0.02061855670103093
$ws.Range("P2").Value = 0.127147766323024
$ws.Range("S2").Value = 0.1134020618556701
$ws.Range("C3").Value = 0.03144654088050314
$ws.Range("J3").Value = 0.01257861635220126
$ws.Range("P3").Value = 0.7232704402515723
$ws.Range("S3").Value = 0.2327044025157233
$ws.Range("J4").Value = 0.02
$ws.Range("P4").Value = 0.78
$ws.Range("S4").Value = 0.2
$ws.Range("B6").Value = 0.06481481481481481
$ws.Range("D6").Value = 0.01388888888888889
$ws.Range("E6").Value = 0.004629629629629629
$ws.Range("F6").Value = 0.05555555555555555
$ws.Range("J6").Value = 0.1666666666666667
$ws.Range("O6").Value = 0.02777777777777778
$ws.Range("Q6").Value = 0.1111111111111111
$ws.Range("R6").Value = 0.1064814814814815
$ws.Range("S6").Value = 0.4490740740740741
$ws.Range("B7").Value = 0.09433962264150944
$ws.Range("D7").Value = 0.0330188679245283
$ws.Range("F7").Value = 0.04245283018867924
$ws.Range("J7").Value = 0.1037735849056604
$ws.Range("O7").Value = 0.004716981132075472
$ws.Range("Q7").Value = 0.2358490566037736
$ws.Range("R7").Value = 0.05660377358490566
$ws.Range("S7").Value = 0.4292452830188679
$ws.Range("B8").Value = 0.08875739644970414
$ws.Range("D8").Value = 0.01972386587771203
$ws.Range("E8").Value = 0.001972386587771203
$ws.Range("F8").Value = 0.07889546351084813
$ws.Range("J8").Value = 0.08284023668639054
$ws.Range("O8").Value = 0.02958579881656805
$ws.Range("Q8").Value = 0.1972386587771203
$ws.Range("R8").Value = 0.08481262327416174
$ws.Range("S8").Value = 0.4161735700197239
$ws.Range("B9").Value = 0.06926406926406926
$ws.Range("D9").Value = 0.01731601731601732
$ws.Range("F9").Value = 0.06060606060606061
$ws.Range("J9").Value = 0.1038961038961039
$ws.Range("O9").Value = 0.01731601731601732
$ws.Range("Q9").Value = 0.1645021645021645
$ws.Range("R9").Value = 0.08225108225108226
$ws.Range("S9").Value = 0.4848484848484849
$ws.Range("B10").Value = 0.1138487680543755
$ws.Range("D10").Value = 0.02293967714528462
$ws.Range("F10").Value = 0.0713678844519966
$ws.Range("J10").Value = 0.1079014443500425
$ws.Range("O10").Value = 0.02888700084961767
$ws.Range("Q10").Value = 0.1945624468988955
$ws.Range("R10").Value = 0.08836023789294817
$ws.Range("S10").Value = 0.3721325403568394
$ws.Range("G11").Value = 0.136986301369863
$ws.Range("J11").Value = 0.07534246575342465
$ws.Range("K11").Value = 0.1883561643835616
$ws.Range("L11").Value = 0.589041095890411
$ws.Range("S11").Value = 0.01027397260273973
$ws.Range("G12").Value = 0.764367816091954
$ws.Range("J12").Value = 0.1781609195402299
$ws.Range("K12").Value = 0.005747126436781609
$ws.Range("L12").Value = 0.02873563218390805
$ws.Range("S12").Value = 0.02298850574712644
$ws.Range("G13").Value = 0.7049180327868853
$ws.Range("J13").Value = 0.2950819672131147
$ws.Range("F15").Value = 0.009259259259259259
$ws.Range("H15").Value = 0.1666666666666667
$ws.Range("I15").Value = 0.07407407407407407
$ws.Range("J15").Value = 0.2962962962962963
$ws.Range("K15").Value = 0.06481481481481481
$ws.Range("M15").Value = 0.009259259259259259
$ws.Range("O15").Value = 0.05092592592592592
$ws.Range("S15").Value = 0.3287037037037037
$ws.Range("F16").Value = 0.02127659574468085
$ws.Range("H16").Value = 0.2127659574468085
$ws.Range("I16").Value = 0.07446808510638298
$ws.Range("J16").Value = 0.4202127659574468
$ws.Range("K16").Value = 0.06914893617021277
$ws.Range("M16").Value = 0.02127659574468085
$ws.Range("O16").Value = 0.04787234042553191
$ws.Range("S16").Value = 0.1329787234042553
$ws.Range("F17").Value = 0.01136363636363636
$ws.Range("H17").Value = 0.2181818181818182
$ws.Range("I17").Value = 0.1022727272727273
$ws.Range("J17").Value = 0.4136363636363636
$ws.Range("K17").Value = 0.07727272727272727
$ws.Range("M17").Value = 0.025
$ws.Range("N17").Value = 0.004545454545454545
$ws.Range("O17").Value = 0.04318181818181818
$ws.Range("S17").Value = 0.1045454545454545
$ws.Range("F18").Value = 0.005050505050505051
$ws.Range("H18").Value = 0.2171717171717172
$ws.Range("I18").Value = 0.101010101010101
$ws.Range("J18").Value = 0.398989898989899
$ws.Range("K18").Value = 0.1060606060606061
$ws.Range("M18").Value = 0.01515151515151515
$ws.Range("O18").Value = 0.06060606060606061
$ws.Range("S18").Value = 0.09595959595959595
$ws.Range("F19").Value = 0.0124031007751938
$ws.Range("H19").Value = 0.2286821705426356
$ws.Range("I19").Value = 0.1077519379844961
$ws.Range("J19").Value = 0.3527131782945737
$ws.Range("K19").Value = 0.1162790697674419
$ws.Range("M19").Value = 0.03178294573643411
$ws.Range("N19").Value = 0.001550387596899225
$ws.Range("O19").Value = 0.06124031007751938
$ws.Range("S19").Value = 0.08759689922480621
